$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.865.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.707.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.55%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.73%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +20.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '675.26'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.86%  '

$ws.Range("E8").Value = '  +7.01%  '

$ws.Range("E9").Value = '  +8.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.703.95'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.35'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.26%  '

$ws.Range("E13").Value = '  +2.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.395.64'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.55%  '

$ws.Range("E16").Value = '  +5.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.619.84'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.88'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +13.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.708.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.551'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '517.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.41%  '

$ws.Range("E24").Value = '  +2.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000210'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.70%  '

$ws.Range("E26").Value = '  +1.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '102.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.68%  '

$ws.Range("E29").Value = '  +13.08%  '

$ws.Range("E30").Value = '  +4.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.27%  '

$ws.Range("E32").Value = '  +0.44%  '

$ws.Range("E33").Value = '  +2.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.31'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.21%  '

$ws.Range("E35").Value = '  +0.35%  '

$ws.Range("E36").Value = '  +7.25%  '

$ws.Range("E37").Value = '  +8.90%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.84'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '617.91'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '42.80'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +27.32%  '

$ws.Range("E41").Value = '  +9.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.970'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.55%  '

$ws.Range("E43").Value = '  +9.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0450'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.427'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +25.65%  '

$ws.Range("E48").Value = '  +3.23%  '

$ws.Range("E49").Value = '  +0.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.64'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.25%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.24%  '
